# Applies the "investigated tree structure and added missingness analysis" edit:
# Adds a second copy of the L0/L1/L2/L3 header table (columns M:P) and a
# missingness/frequency analysis column (N4:N14) next to the existing
# Node/Parent/L0/L1/L2/L3 tables on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 3), columns M:P, reusing existing strings L0,L1,L2,L3 ---
$ws.Range("M3").Value = "L0"
$ws.Range("N3").Value = "L1"
$ws.Range("O3").Value = "L2"
$ws.Range("P3").Value = "L3"

# --- New data: M4:M4 and N4:N14 (missingness / value-count analysis) ---
$ws.Range("M4").Value = 2833
$ws.Range("N4").Value = 2
$ws.Range("N5").Value = 106
$ws.Range("N6").Value = 150
$ws.Range("N7").Value = 191
$ws.Range("N8").Value = 206
$ws.Range("N9").Value = 220
$ws.Range("N10").Value = 225
$ws.Range("N11").Value = 242
$ws.Range("N12").Value = 830
$ws.Range("N13").Value = 839
$ws.Range("N14").Value = 2557

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("O4").Select()

# --- Window position tweak recorded alongside the edit ---
$win = $excel.ActiveWindow
$win.Left = 28680
$win.Top = -120
